$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top; this shifts the existing
# rows 1-16 down to rows 3-18 (all data + formatting moves with them).
$ws.Rows("1:2").Insert()

# Row 3 is the old row 1 (the text header row) and already kept its bold/
# bordered "header" formatting through the insert. Give new row 1 the same
# look by copying that formatting across before filling in its values.
$ws.Range("A3:N3").Copy()
$ws.Range("A1:N1").PasteSpecial(-4122)

# New row 1: numeric column-index header (0-13).
$headerValues = @(0, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13)
for ($i = 0; $i -lt $headerValues.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headerValues[$i]
}

# New row 2: blank "spacer" row, only E2 holds "Washer".
for ($col = 1; $col -le 14; $col++) {
    $ws.Cells.Item(2, $col).Value = $null
}
$ws.Range("E2").Value = "Washer"

# Row 3 (the old header row) loses its M/N labels - they become blank.
$ws.Range("M3").Value = $null
$ws.Range("N3").Value = $null

# K3 (old K1) was already an empty placeholder cell; keep it a true blank
# instead of the empty-string artifact row-insert can leave behind.
$ws.Range("K3").Value = $null
